$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (style s="2", date format) from A226 to the new date cells A227:A229
$ws.Range("A226").Copy() | Out-Null
$ws.Range("A227:A229").PasteSpecial(-4122) | Out-Null

# Row 227
$ws.Cells.Item(227, 1).Value = 44301
$ws.Cells.Item(227, 2).Value = 10
$ws.Cells.Item(227, 3).Value = 10
$ws.Cells.Item(227, 4).Value = 8
$ws.Cells.Item(227, 5).Value = 126
$ws.Cells.Item(227, 6).Value = 51
$ws.Cells.Item(227, 7).Value = 33
$ws.Cells.Item(227, 8).Value = 24
$ws.Cells.Item(227, 9).Value = 9
$ws.Cells.Item(227, 10).Value = 9
$ws.Cells.Item(227, 11).Value = 1
$ws.Cells.Item(227, 12).Value = 31
$ws.Cells.Item(227, 13).Value = 36
$ws.Cells.Item(227, 14).Value = 56
$ws.Cells.Item(227, 15).Value = 2
$ws.Cells.Item(227, 16).Value = 7
$ws.Cells.Item(227, 17).Value = 5
$ws.Cells.Item(227, 18).Value = 40
$ws.Cells.Item(227, 19).Value = 1
$ws.Cells.Item(227, 20).Value = 7
$ws.Cells.Item(227, 21).Value = 38
$ws.Cells.Item(227, 22).Value = 333
$ws.Cells.Item(227, 23).Value = 7
$ws.Cells.Item(227, 24).Value = 24
$ws.Cells.Item(227, 25).Value = 25
$ws.Cells.Item(227, 26).Value = 44
$ws.Cells.Item(227, 27).Value = 4
$ws.Cells.Item(227, 28).Value = 11
$ws.Cells.Item(227, 29).Value = 7
$ws.Cells.Item(227, 30).Value = 7
$ws.Cells.Item(227, 31).Value = 38
$ws.Cells.Item(227, 32).Value = 8
$ws.Cells.Item(227, 33).Value = 113
$ws.Cells.Item(227, 34).Value = 11
$ws.Cells.Item(227, 35).Value = 25
$ws.Cells.Item(227, 36).Value = 8
$ws.Cells.Item(227, 37).Value = 32
$ws.Cells.Item(227, 38).Value = 25
$ws.Cells.Item(227, 39).Value = 60
$ws.Cells.Item(227, 40).Value = 2
$ws.Cells.Item(227, 41).Value = 45
$ws.Cells.Item(227, 42).Value = 1366
$ws.Cells.Item(227, 43).Value = 20
$ws.Cells.Item(227, 44).Value = 5
$ws.Cells.Item(227, 45).Value = 2
$ws.Cells.Item(227, 46).Value = 1
$ws.Cells.Item(227, 47).Value = 1
$ws.Cells.Item(227, 48).Value = 4
$ws.Cells.Item(227, 49).Value = 1
$ws.Cells.Item(227, 50).Value = 0

# Row 228
$ws.Cells.Item(228, 1).Value = 44302
$ws.Cells.Item(228, 2).Value = 12
$ws.Cells.Item(228, 3).Value = 11
$ws.Cells.Item(228, 4).Value = 7
$ws.Cells.Item(228, 5).Value = 128
$ws.Cells.Item(228, 6).Value = 57
$ws.Cells.Item(228, 7).Value = 32
$ws.Cells.Item(228, 8).Value = 30
$ws.Cells.Item(228, 9).Value = 10
$ws.Cells.Item(228, 10).Value = 12
$ws.Cells.Item(228, 11).Value = 1
$ws.Cells.Item(228, 12).Value = 24
$ws.Cells.Item(228, 13).Value = 39
$ws.Cells.Item(228, 14).Value = 47
$ws.Cells.Item(228, 15).Value = 2
$ws.Cells.Item(228, 16).Value = 6
$ws.Cells.Item(228, 17).Value = 4
$ws.Cells.Item(228, 18).Value = 29
$ws.Cells.Item(228, 19).Value = 3
$ws.Cells.Item(228, 20).Value = 10
$ws.Cells.Item(228, 21).Value = 42
$ws.Cells.Item(228, 22).Value = 266
$ws.Cells.Item(228, 23).Value = 12
$ws.Cells.Item(228, 24).Value = 19
$ws.Cells.Item(228, 25).Value = 19
$ws.Cells.Item(228, 26).Value = 44
$ws.Cells.Item(228, 27).Value = 4
$ws.Cells.Item(228, 28).Value = 8
$ws.Cells.Item(228, 29).Value = 12
$ws.Cells.Item(228, 30).Value = 8
$ws.Cells.Item(228, 31).Value = 46
$ws.Cells.Item(228, 32).Value = 7
$ws.Cells.Item(228, 33).Value = 106
$ws.Cells.Item(228, 34).Value = 14
$ws.Cells.Item(228, 35).Value = 22
$ws.Cells.Item(228, 36).Value = 7
$ws.Cells.Item(228, 37).Value = 34
$ws.Cells.Item(228, 38).Value = 27
$ws.Cells.Item(228, 39).Value = 55
$ws.Cells.Item(228, 40).Value = 3
$ws.Cells.Item(228, 41).Value = 35
$ws.Cells.Item(228, 42).Value = 1282
$ws.Cells.Item(228, 43).Value = 20
$ws.Cells.Item(228, 44).Value = 2
$ws.Cells.Item(228, 45).Value = 2
$ws.Cells.Item(228, 46).Value = 1
$ws.Cells.Item(228, 47).Value = 1
$ws.Cells.Item(228, 48).Value = 3
$ws.Cells.Item(228, 49).Value = 1
$ws.Cells.Item(228, 50).Value = 0

# Row 229
$ws.Cells.Item(229, 1).Value = 44303
$ws.Cells.Item(229, 2).Value = 13
$ws.Cells.Item(229, 3).Value = 6
$ws.Cells.Item(229, 4).Value = 5
$ws.Cells.Item(229, 5).Value = 112
$ws.Cells.Item(229, 6).Value = 60
$ws.Cells.Item(229, 7).Value = 28
$ws.Cells.Item(229, 8).Value = 28
$ws.Cells.Item(229, 9).Value = 10
$ws.Cells.Item(229, 10).Value = 12
$ws.Cells.Item(229, 11).Value = 1
$ws.Cells.Item(229, 12).Value = 21
$ws.Cells.Item(229, 13).Value = 43
$ws.Cells.Item(229, 14).Value = 48
$ws.Cells.Item(229, 15).Value = 2
$ws.Cells.Item(229, 16).Value = 3
$ws.Cells.Item(229, 17).Value = 3
$ws.Cells.Item(229, 18).Value = 24
$ws.Cells.Item(229, 19).Value = 2
$ws.Cells.Item(229, 20).Value = 8
$ws.Cells.Item(229, 21).Value = 46
$ws.Cells.Item(229, 22).Value = 243
$ws.Cells.Item(229, 23).Value = 11
$ws.Cells.Item(229, 24).Value = 17
$ws.Cells.Item(229, 25).Value = 16
$ws.Cells.Item(229, 26).Value = 39
$ws.Cells.Item(229, 27).Value = 4
$ws.Cells.Item(229, 28).Value = 6
$ws.Cells.Item(229, 29).Value = 13
$ws.Cells.Item(229, 30).Value = 5
$ws.Cells.Item(229, 31).Value = 41
$ws.Cells.Item(229, 32).Value = 8
$ws.Cells.Item(229, 33).Value = 124
$ws.Cells.Item(229, 34).Value = 9
$ws.Cells.Item(229, 35).Value = 24
$ws.Cells.Item(229, 36).Value = 4
$ws.Cells.Item(229, 37).Value = 32
$ws.Cells.Item(229, 38).Value = 22
$ws.Cells.Item(229, 39).Value = 51
$ws.Cells.Item(229, 40).Value = 3
$ws.Cells.Item(229, 41).Value = 36
$ws.Cells.Item(229, 42).Value = 1207
$ws.Cells.Item(229, 43).Value = 18
$ws.Cells.Item(229, 44).Value = 2
$ws.Cells.Item(229, 45).Value = 1
$ws.Cells.Item(229, 46).Value = 1
$ws.Cells.Item(229, 47).Value = 1
$ws.Cells.Item(229, 48).Value = 2
$ws.Cells.Item(229, 49).Value = 1
$ws.Cells.Item(229, 50).Value = 0
